# Updates rows 2-19: columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
# The edit shuffles the per-row data (dates + associated volume/price figures)
# between rows while leaving all other columns (market/product metadata) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(Fecha(serial), Volumen, PrecioMin, PrecioMax, PrecioPromPond, PrecioKg)
# Row 8 is intentionally omitted: it is unchanged by this edit.
$data = @{
    2  = @(44455, 160, 13000, 14000, 13500, 6750)
    3  = @(44494, 200, 11500, 12000, 11750, 5875)
    4  = @(44466, 160, 13500, 14000, 13750, 6875)
    5  = @(44462, 140, 13000, 14000, 13500, 6750)
    6  = @(44497, 400, 11500, 12000, 11750, 5875)
    7  = @(44445, 160, 14000, 15000, 14500, 7250)
    9  = @(44495, 300, 11000, 12000, 11500, 5750)
    10 = @(44468, 300, 13000, 14000, 13500, 6750)
    11 = @(44498, 240, 11000, 11500, 11250, 5625)
    12 = @(44489, 400, 11500, 12000, 11750, 5875)
    13 = @(44490, 160, 11500, 12000, 11750, 5875)
    14 = @(44491, 200, 11500, 12000, 11750, 5875)
    15 = @(44446, 300, 14000, 15000, 14500, 7250)
    16 = @(44459, 160, 13000, 14000, 13500, 6750)
    17 = @(44452, 200, 13000, 14000, 13500, 6750)
    18 = @(44463, 100, 13000, 14000, 13500, 6750)
    19 = @(44454, 300, 13000, 14000, 13500, 6750)
}

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $fecha = $epoch.AddDays([double]$vals[0])
    $ws.Cells.Item($row, 4).Value = $fecha

    $ws.Cells.Item($row, 13).Value = $vals[1]
    $ws.Cells.Item($row, 14).Value = $vals[2]
    $ws.Cells.Item($row, 15).Value = $vals[3]
    $ws.Cells.Item($row, 16).Value = $vals[4]
    $ws.Cells.Item($row, 19).Value = $vals[5]
}
